$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$t  = $ws.ListObjects.Item(1)

# ------------------------------------------------------------------
# 1. Add two new rows to the table (rows 14 and 15) BEFORE we touch
#    any text, so that later writes land on real worksheet rows and
#    the table/ worksheet dimension grows accordingly.
# ------------------------------------------------------------------
$t.ListRows.Add() | Out-Null
$t.ListRows.Add() | Out-Null

# ------------------------------------------------------------------
# 2. Fill in row 14 : "Positive EV (4%)"
# ------------------------------------------------------------------
$ws.Range("E14").Value = "Positive EV (4%)"
$ws.Range("F14").Formula = "=G14/(I14*J14)"
$ws.Range("G14").Value = -26.13
$ws.Range("H14").Formula = "=Table1[[#This Row],[Matches]]*Table1[[#This Row],[Cost per match]]"
$ws.Range("I14").Value = 925
$ws.Range("J14").Value = 1

# ------------------------------------------------------------------
# 3. Fix the misspelled "Aribitrage" entry (row 2) and refresh its
#    underlying numbers to the new arbitrage result.
# ------------------------------------------------------------------
$ws.Range("E2").Value = "Arbitrage"
$ws.Range("I2").Value = 497
$ws.Range("J2").Value = 1
$ws.Range("G2").Formula = "=3.453829"
$ws.Range("F2").Style = "Percent"

# ------------------------------------------------------------------
# 4. Fill in row 15 : "Arbitrage (Compounded)"
# ------------------------------------------------------------------
$ws.Range("E15").Value = "Arbitrage (Compounded)"
$ws.Range("F15").Formula = "=G15/(I15*J15)"
$ws.Range("G15").Value = 3.2351424949999998
$ws.Range("H15").Formula = "=Table1[[#This Row],[Matches]]*Table1[[#This Row],[Cost per match]]"
$ws.Range("I15").Value = 50
$ws.Range("J15").Value = 1

# ------------------------------------------------------------------
# 5. Column width adjustments (Method / ROI / Profit columns).
# ------------------------------------------------------------------
$ws.Columns("E").ColumnWidth = 27.17
$ws.Columns("F").ColumnWidth = 22.5
$ws.Columns("G").ColumnWidth = 11.33

# ------------------------------------------------------------------
# 6. Keep the active selection pointing at the new last row, matching
#    where the author's cursor ended up.
# ------------------------------------------------------------------
$ws.Range("E15").Select() | Out-Null
